$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update the metadata blob stored in the A1 cell comment (produced by the
#    CEIC Data Manager add-in). The text is a base64-encoded, gzip-compressed
#    XML payload; only its content (sheet name / timestamps / ids baked in by
#    the add-in) changed between revisions, so we replace it verbatim.
$newCommentText = 'gjQAAB+LCAAAAAAAAAPtW+tv29YV/1cIfeqAyCRlObGNaxayZHvaLNu1lDnpl4Air20uFKnyYVvfEmBrnKbFFmQp0Gc27NGtW7sU7ZA2dtr/pbDk5FP+hZ374FOULXoJsBQOjJj3PO4999zz+PGGQa/vdUxhBzuuYVtzBXlCKgjY0mzdsLbmCr63WZQvFl5X0MKehs011VE72ANhAbQsd3bPNeYK257XnRXF3d3did3JCdvZEkuSJItXGstNbRt31KJhuZ5qabgQaumnaxUUVNU7DeypuuqpTHOuUG/WJ6rY0GpAa6iWuoWdiXnfNSzsuguWZ3gGdommg1UPV2uNX7GNKaWJixMyEofokeS8b5g6k0tIMjqXg2Vxy+hgpSTJ00VppihNtmRpdmp6VipNlC9OvxkohoJoWXW9JnZ2DI0Smp7a6VJ1aUaalKWpaamExEwhmCtygIJWTX0d7xgu1qvYNN1cHhH5AVY0D3adz5kSEmO6fKKzm7DkqN3tluGZOJ8Z6415oW1dEMjvDvyGaCkX4a+pOVkKTIzmVtCi7WAN3HomS1fw7qrDvd3qLgO3tW04Xq+m9nLPddnFzmqX+C6fqoJqtuVVTOx4l7sQAliHCAGG4jk+RuIIZqRUM1wNng3Lx7qyqZpuXCnBRBu2c93tqhpegfQWyRy7lmmrOsShZ7ieoUWLDjHQmmN3YUZYfN429UWYlQtnMMKZ6xa4mCw7b9vXI+uymIieKj1fONOO6gXiQ3TU3LZ3Vy2z1/TbruYYbazX5gPpTB4iecq1q77r2R2wIiIhRotRdF1sNMQe/IHsTDNRDWtGRzXXTHClq0zCXAkCqvievWl4Vdv0O5YbWJaiog3YVwvvhfsMx2gVjtgirretuhXIM2dnspIK6/ZuuOYwg7oiRq64WnDow4y0cA1owSEOc+i5kF0uGiY0j/iJxKjJ2GhuY+xlBgbjIFInF0k7UuZ7K36nDUnWhkzboau6SIz4CKIVIh7sUiRoMUX605KkWfoDdoRstGDpo+UCJoLlYmsp8hQSUyQEezLnTdW6DtQNw9teqQR7yeAg5oGR8sM8BPnbNdUeJYdeitNQ3dJMX8esLNStTRqixDZ2qCPZaIi0DJmuINXqtXpdKNquMevBw1wB2vis6zkAFAqKZvuW5/RI/UAiFz1Nx/XbFl1ANcfW2XTwWz7gk96ib2lVWx9/NZ1557JleONbaPsOK4rjq1DvkfrouzVMKg0t/WPra3n25Dq5xDsW7tiWoY3vbXAysV4/w0bcIKvG1sAsv8aWN6G7s+ZHcn1sNQfAJbS7XMtUXNfWDBqsPD30mL44ImVqeFP1TQB2HjTarbD2psmo4l5Py8RJ6LJjBhVQIbDZBdys6Z0JDSAEwYYTmt0hBBHg6kYTiXF5goM0vGBtLavWlg9II6wraXpYf0mXbDmq5ZLthMAiVYqzhVBQpxjgUVjxWvVpILDiZQMXiSk51MKdru2oZgMcYyzysOOoCQBJQ/W2+Qh6m4m1wMlipBpqJS0LDD9NjDYptg2S8LxMpohUiOyFIfRIJqIhsssGpKVZVU2j7bCqGrTyLB4cWAQRg/pLNpcTLgZnAG9p0H1/iXsEukcDTqchKwcMFsCkkCrN9fJ0qTwzOQPIhowR3fGSY7uuULM7mMA9ATCd7mue8NpSbe1n0PaoU3jjGUM+LokWg1pOl3sDbEpSkgKAKbYMaBXDgiEnUlDe8KH6YMfsxWTZjpZtDUT7Hz3pv/3b4z88OHr82eDGQf/G4fPDfbDx+eFttismhlpq28TUmtb89LQ0WYZACkmIeE+kCJhsktKuXqXANxwj/lpHB9WFenVpeZ4WjJAYqLOeIZI3xp7tR8MmM58uRM9MDI6aiSitoADxcYIb60EKeYPbwUnpOH+UIvPF8cFnxwdfjNTmDovAlDwzQ96PT8Va8C4sF+XSCKy1nCjyRLhclKaKpVJMOCWD1lmFD/1U1xV4056R4UVRDou1HsZtllCaxWdqqVtiSo+Rqgz/hCEQHwdMGvUtyIuQzfIgNmBOHPznztN/309Ice9ySnIWMI7iFLKYGAzo1CvrLaG5enm9uiC0FpokTiJeTI5NfoIwXz3MpHjQhJl2QYAGDg1LKDRU54LwC9/CF4Qm7l4Q4C2oEMvDRMhlUdly/9PEaYuhLPldVmdiChE1QzIqOFkaGUWH8qhvs6pPxM3Q4EXpiz9nKfC91CIAyy9DwjClNJTgMFKMz5d47+bR4x+OHj8ePLzb//ZmYga+TngZAGEPyRUfhlkAFZD3lxQFbTSpP69L12L9hBPJi9OabVieq0yV6DsTHyFQlcls9Deqd6DF0Ympv4CeoqCfq+7CnsfzXFlBYpIAdnZV6K529K4ZElhJj/z67JNPBx99M3j/4dNbn/ff+Xv/3fePDx48/fJPLAkH9x8O3vuSF/10X6C2kDdYBvoEeiuiCSQ5BdKrhR9v3BMs2xMAYgg+LVA/3vggNhkxlIKRaGaAcKEhSROGROPKRE+ImRLakNALVVjDr5KONhlK8J5mdw0tWuTNIpmKZB9lvFZvFX0XCzagJ+jnKeFIeVw9rsI67NolqSSXOJdZQ7bQVt2Y65dMu62aQsCgFw4pkYTWyQqRLF1vaXl1vrIciTAjVh0dOyQM2QMKICTpMHU3GAWhFqMAF4Ce5pvkjmhIbJgVzhyrZCK/bNms6KQIZl9NJCRQ1XccBo4sfrHf9LsAfoOLudF8elkZw7srDJvGEXA0rteSfBjHuNAXk2xCoHxamjiLlam6S+51GHxdIa6JhsBLXHCCO/jdPQNeOwYciUjqzoLj2E5m8Yk4gVgDkDNUFDHyeChDz5ShbD06q4AQFLwXg7anpPKlBNperNRnhZ6nzwp1aHmwNx8iN4D6CWLVJperKgN0wzD8jBOdgM8bY8LztFwWOm/Ylrc9EpuD7c8PP+x/901//5Ont74+OngHhkfffn705F7/0V/7998e7D8CyvGdfx3/807/9/tHTz4+/sv357D9xcB2+u8qY1yRMiRefvmwXZounw7budA5bB8Xtg+n4Jiw/ETF0bC7kRt2Z2qcCLuHTRsDdA8+3s8BujsjQXdnNOg+/uD7o29vvSjQLY8G3fK1WFsZBt1yqXyOus9R9yuFuuVz1P2qoO5U9fk/Rt3TyTvujTUAy2uOAWWsJ1QcqCYmBpDcALDsgCnwhD3VME1DA+TsQJ6vOpmg+2zznIC5Ny6OCbqHBLNQ9wbG10eC7sHt7/q/u3t08Gjwx8PBu7eg/D0/3IcNPT+8TdH3wdO/3QTo3b93k0DvBz8c3/4HPDy7dbf/1afHH/6GPN+7SegPvj5H4i8MiU9RhF06FYmT2/OiXH7ZF+iTU2NcoDOhcyQ+LhJnaSm0MSSsRVq36glNiA1HV88Az/PPNhqzb+TG7JkaJ2L2oaI0BmR/duOrMSF76iPBYeSeEhgB4EdPkx++l0bD99K1WH8ahu+TU1Pn8P0cvr9S8L30E4DvfO2fPH5PlZ+Xgt/F6OuWsDLaNWxiL9/n4GKk3bB3zqwLh59Xte6umjp3Zr4PdEK3RBPEv4kngfKSPolnMVhxHMBc5OvZ3N+wB199ravWVk4j2Q6pIvl8ClbnH5kvGo7rXSGNij8xytWQcpWB1yvkQ2r2QMdXldI0x6wgIcanFxN2Bjntsf+xYJvLRsfI+VGVFOR9chJwZrfL4F09XwSRnrOC9wB8xmaAatn+NfQT9hVintlYIEORDfXJl7+usbXt5TXsUlvFOm5LRa2NS8WyLk0XZzCeLMoy/K1qJRJY5LthPjlUFAPv5lxEDA4s+p9Eyn8B90dTy4I0AAA='
$ws.Range("A1").Comment.Text($newCommentText)

# 3. Change the custom number format used by B27:D281 from "0.000" to "###0.000"
$ws.Range("B27:D281").NumberFormat = "###0.000"

# 4. Update the label in A11
$ws.Range("A11").Value = "Function Information"
